$wb = $excel.ActiveWorkbook

# --- Sheet "Item" ---
$item = $wb.Worksheets.Item("Item")
$item.Rows.Item(1).Insert()
$item.Range("A1").Value = "# Dummy Row (Ignored by ExcelBinder)"

# --- Sheet "Skill" ---
$skill = $wb.Worksheets.Item("Skill")
$skill.Rows.Item(1).Insert()
$skill.Range("A1").Value = "# Dummy Row (Ignored by ExcelBinder)"
